$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet2 ("MC3PO Sample"): insert 3 new columns (Exposures, Exposure (s),
# Total Exp (m)) before the old "Flat" column (G), pushing Flat..Comment from
# columns G..L to J..O.
$ws2.Range("G1:I1").EntireColumn.Insert()

# Row numbers/dates for the new MC3PO rows (dates don't touch the shared
# string table, so these can be written in any order).
$ws2.Range("A11").Value = 39307
$ws2.Range("A12").Value = 36744
$ws2.Range("A13").Value = 36744
$ws2.Range("A14").Value = 37829
$ws2.Range("A15").Value = 37889
$ws2.Range("A16").Value = 39219
$ws2.Range("A17").Value = 39281

# Exposure (s) values -- plain numbers, order-independent.
$ws2.Range("H2").Value = 60
$ws2.Range("H3").Value = 240
$ws2.Range("H4").Value = 180
$ws2.Range("H5").Value = 240
$ws2.Range("H6").Value = 180
$ws2.Range("H7").Value = 3
$ws2.Range("H8").Value = 3
$ws2.Range("H9").Value = 3
$ws2.Range("H10").Value = 3
$ws2.Range("H11").Value = 180
$ws2.Range("H12").Value = 480
$ws2.Range("H13").Value = 180
$ws2.Range("H14").Value = 240
$ws2.Range("H15").Value = 360
$ws2.Range("H16").Value = 150
$ws2.Range("H17").Value = 180

# Cluster bands / flags for the new rows -- reuse of already-existing shared
# strings ("W-C-RC", "W-C-IC", "W-J-V", "W-J-B", "W-S-Z+", "wd") so order
# relative to brand-new strings below does not matter.
$ws2.Range("C11").Value = "W-J-B"
$ws2.Range("D11").Value = "wd"
$ws2.Range("C12").Value = "W-C-RC"
$ws2.Range("D12").Value = "wd"
$ws2.Range("C13").Value = "W-C-IC"
$ws2.Range("D13").Value = "wd"
$ws2.Range("C14").Value = "W-J-B"
$ws2.Range("D14").Value = "wd"
$ws2.Range("C15").Value = "W-J-V"
$ws2.Range("D15").Value = "wd"
$ws2.Range("C16").Value = "W-C-RC"
$ws2.Range("D16").Value = "wd"
$ws2.Range("C17").Value = "W-S-Z+"
$ws2.Range("D17").Value = "wd"

# --- Text values that mint brand-new shared strings: written in the same
# order the original author entered them so the shared string table lines up.
$ws2.Range("E3").Value = "SUPA00232340"
$ws2.Range("F3").Value = "SUPA00232429"
$ws2.Range("E4").Value = "SUPA00330440"
$ws2.Range("F4").Value = "SUPA00330529"
$ws2.Range("H1").Value = "Exposure (s)"
$ws2.Range("E5").Value = "SUPA00487380"
$ws2.Range("F5").Value = "SUPA00487459"
$ws2.Range("E6").Value = "SUPA00557410"
$ws2.Range("F6").Value = "SUPA00557509"
$ws2.Range("E7").Value = "SUPA00561680"
$ws2.Range("F7").Value = "SUPA00561689"
$ws2.Range("E8").Value = "SUPA00561700"
$ws2.Range("F8").Value = "SUPA00561709"
$ws2.Range("E9").Value = "SUPA00561720"
$ws2.Range("F9").Value = "SUPA00561729"
$ws2.Range("E10").Value = "SUPA00561740"
$ws2.Range("F10").Value = "SUPA00561749"
$ws2.Range("E11").Value = "SUPA00562030"
$ws2.Range("F11").Value = "SUPA00562109"
$ws2.Range("B4").Value = "RXC J2228.6+2037"
$ws2.Range("B5").Value = "RXC J2228.6+2038"
$ws2.Range("B6").Value = "RXC J2228.6+2039"
$ws2.Range("B7").Value = "RXC J2228.6+2040"
$ws2.Range("B8").Value = "RXC J2228.6+2041"
$ws2.Range("B9").Value = "RXC J2228.6+2042"
$ws2.Range("B10").Value = "RXC J2228.6+2043"
$ws2.Range("B11").Value = "RXC J2228.6+2044"
$ws2.Range("E2").Value = "SUPA00395940"
$ws2.Range("F2").Value = "SUPA00396029"
$ws2.Range("E12").Value = "SUPA00022752"
$ws2.Range("F12").Value = "SUPA00022809"
$ws2.Range("E13").Value = "SUPA00022852"
$ws2.Range("F13").Value = "SUPA00022909"
$ws2.Range("E14").Value = "SUPA00232540"
$ws2.Range("F14").Value = "SUPA00232559"
$ws2.Range("E15").Value = "SUPA00242810"
$ws2.Range("F15").Value = "SUPA00242869"
$ws2.Range("E16").Value = "SUPA00543720"
$ws2.Range("G1").Value = "Exposures"
$ws2.Range("F16").Value = "SUPA00543769"
$ws2.Range("E17").Value = "SUPA00557520"
$ws2.Range("F17").Value = "SUPA00557609"
$ws2.Range("B13").Value = "MACS J2243.3-0936"
$ws2.Range("B14").Value = "MACS J2243.3-0937"
$ws2.Range("B15").Value = "MACS J2243.3-0938"
$ws2.Range("B16").Value = "MACS J2243.3-0939"
$ws2.Range("B17").Value = "MACS J2243.3-0940"
$ws2.Range("I1").Value = "Total Exp (m)"

# B12 reuses the cluster name that used to live at B28 further down the
# (mostly empty) sheet; move it up and clear the stale occurrence.
$ws2.Range("B12").Value = "MACS J2243.3-0935"
$ws2.Range("B28").ClearContents()

# Exposures / Total Exp (m) formulas, filled as one contiguous block so Excel
# records a single shared-formula group (G3:G17 / I3:I17), matching how the
# author filled the column down.
$ws2.Range("G3:G17").Formula = "=(RIGHT(F3,LEN(F3)-4)-RIGHT(E3,LEN(E3)-4)+1)/10"
$ws2.Range("I3:I17").Formula = "=G3*H3/60"

# Row 2 uses its own (non-shared) formulas.
$ws2.Range("G2").Formula = "=(RIGHT(F2,LEN(F2)-4)-RIGHT(E2,LEN(E2)-4)+1)/10"
$ws2.Range("I2").Formula = "=G2*H2/60"

# Rows 12/13 were hand-entered (6 exposures) instead of left as formulas.
$ws2.Range("G12").Value = 6
$ws2.Range("G13").Value = 6

# Selection moved to B2 the last time the sheet was saved.
$ws2.Range("B2").Select()
